$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (the oldest-quarter-first table grows a newer
# column on the left); existing D:K data shifts right to F:M.
$ws.Range("D:E").Insert()

# The new D:E columns should carry the same number/date formatting as the data that used
# to live there and now lives in F:G, not the plain format Excel copies from column C.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Insert() stamped the (format-less) whole-column style onto header rows 5 and 6, which
# have no data columns at all in this table - strip that back off.
$ws.Range("D5:E6").Clear()

# Populate the two new columns with the newest two quarters of data.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1635000
$ws.Range("E8").Value = 1872400
$ws.Range("D9").Value = 1392000
$ws.Range("E9").Value = 1551700
$ws.Range("D10").Value = 243000
$ws.Range("E10").Value = 320700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 7200
$ws.Range("E14").Value = -106300
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 1508200
$ws.Range("E17").Value = 1557900
$ws.Range("D18").Value = 126800
$ws.Range("E18").Value = 314500
$ws.Range("D20").Value = 5300
$ws.Range("E20").Value = 5500
$ws.Range("D21").Value = 282500
$ws.Range("E21").Value = 473600
$ws.Range("D22").Value = 58000
$ws.Range("E22").Value = 58600
$ws.Range("D23").Value = 74100
$ws.Range("E23").Value = 261400
$ws.Range("D24").Value = 20800
$ws.Range("E24").Value = 65900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 53300
$ws.Range("E26").Value = 195500
$ws.Range("D27").Value = 53300
$ws.Range("E27").Value = 195500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = -400
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -5300
$ws.Range("E32").Value = -5500
$ws.Range("D33").Value = 53300
$ws.Range("E33").Value = 195100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 53300
$ws.Range("E35").Value = 195100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 178800
$ws.Range("E41").Value = 156700
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 782200
$ws.Range("E43").Value = 1022000
$ws.Range("D44").Value = 711400
$ws.Range("E44").Value = 724400
$ws.Range("D45").Value = 35000
$ws.Range("E45").Value = 35200
$ws.Range("D46").Value = 1707400
$ws.Range("E46").Value = 1938300
$ws.Range("D47").Value = 8800
$ws.Range("E47").Value = 13200
$ws.Range("D48").Value = 3482100
$ws.Range("E48").Value = 3456700
$ws.Range("D49").Value = 2631200
$ws.Range("E49").Value = 2647900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1167900
$ws.Range("E52").Value = 1171200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 8997400
$ws.Range("E54").Value = 9227300
$ws.Range("D57").Value = 636500
$ws.Range("E57").Value = 709400
$ws.Range("D58").Value = 125900
$ws.Range("E58").Value = 900
$ws.Range("D59").Value = 355900
$ws.Range("E59").Value = 363400
$ws.Range("D60").Value = 1118300
$ws.Range("E60").Value = 1073700
$ws.Range("D61").Value = 3104400
$ws.Range("E61").Value = 3336400
$ws.Range("D62").Value = 1942500
$ws.Range("E62").Value = 1894300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 6165200
$ws.Range("E66").Value = 6304400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1070500
$ws.Range("E72").Value = 1050500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 2832200
$ws.Range("E76").Value = 2922900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 53300
$ws.Range("E81").Value = 195100
$ws.Range("D83").Value = 150400
$ws.Range("E83").Value = 153600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 321500
$ws.Range("E89").Value = 322000
$ws.Range("D91").Value = -110700
$ws.Range("E91").Value = -98500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -110700
$ws.Range("E94").Value = -95700
$ws.Range("D96").Value = -33300
$ws.Range("E96").Value = -33400
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -188700
$ws.Range("E100").Value = -213700
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = -100
$ws.Range("D102").Value = 22100
$ws.Range("E102").Value = 12500
